$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Google-Sheet export added a new metadata row ("skos:prefLabel" / "WEAVE")
# right before the existing "dct:description" row, pushing every row from the
# old row 9 onward down by one (old row 53 -> new row 54).
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "skos:prefLabel"
$ws.Range("B9").Value = "WEAVE"
$ws.Range("C9").Value = "prefLabel of controlled vocabulary"
